$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(4).Insert()
Write-Host "Inserted"
Write-Host $ws.Cells.Item(3,4).Value()
Write-Host $ws.Cells.Item(3,5).Value()
